# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting the newly scraped numbers from the data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 67
$ws.Range("F5").Value = 243
$ws.Range("F10").Value = 49
$ws.Range("F13").Value = 2280
$ws.Range("F15").Value = 34
$ws.Range("F17").Value = 527
$ws.Range("F18").Value = 160
$ws.Range("F19").Value = 82
$ws.Range("F20").Value = 41
$ws.Range("F22").Value = 1759
$ws.Range("F23").Value = 3898
$ws.Range("F25").Value = 63
$ws.Range("F27").Value = 1163
$ws.Range("F28").Value = 219
$ws.Range("F29").Value = 2060
$ws.Range("F32").Value = 92
$ws.Range("F33").Value = 284
$ws.Range("F35").Value = 459
$ws.Range("F36").Value = 680
$ws.Range("F38").Value = 404

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 25

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 67
$ws.Range("F5").Value = 243
$ws.Range("F10").Value = 49
$ws.Range("F13").Value = 2280
$ws.Range("F15").Value = 25
$ws.Range("F16").Value = 34
$ws.Range("F18").Value = 527
$ws.Range("F19").Value = 160
$ws.Range("F20").Value = 82
$ws.Range("F21").Value = 41
$ws.Range("F23").Value = 1759
$ws.Range("F24").Value = 3898
$ws.Range("F26").Value = 63
$ws.Range("F28").Value = 1163
$ws.Range("F29").Value = 219
$ws.Range("F30").Value = 2060
$ws.Range("F33").Value = 92
$ws.Range("F34").Value = 284
$ws.Range("F36").Value = 459
$ws.Range("F37").Value = 680
$ws.Range("F39").Value = 404
